# v0.2.0 version of the CBS IG
$wb = $excel.ActiveWorkbook

# --- 1. Bump Version + Date on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.2.0"
$meta.Range("B8").Value = "2022-02-08T17:12:45-05:00"

# --- 2. Add the new "Include from Case Based Surve" sheet after "Include ValueSets 2" ---
$srcInclude = $wb.Worksheets.Item("Include ValueSets")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Include from Case Based Surve"

# Match the header & body cell styles used by the other "Include ValueSets" sheets
$srcInclude.Range("A1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)  # xlPasteFormats

$srcInclude.Range("A2").Copy()
$newSheet.Range("A2:B4").PasteSpecial(-4122)  # xlPasteFormats

# Match the column widths used by the other "Include ValueSets" sheets (30.703125 / 50.703125
# characters -- closest values reachable through ColumnWidth's internal 1/6-character rounding)
$newSheet.Columns.Item(1).ColumnWidth = 29.75
$newSheet.Columns.Item(2).ColumnWidth = 49.75

# --- 3. Populate the new sheet's data grid ---
$newSheet.Range("A1").Value = "Concept"
$newSheet.Range("B1").Value = "Description"
$newSheet.Range("A2").Value = "lab-interpretative-report"
$newSheet.Range("B2").Value = "Lab Interpretive Report"
$newSheet.Range("A3").Value = ""
$newSheet.Range("B3").Value = ""
$newSheet.Range("A4").Value = "System URI"
$newSheet.Range("B4").Value = "http://cbsig.chai.gatech.edu/CodeSystem/cbs-temp-code-system"

# --- 4. Restore original active sheet/tab selection (Metadata) ---
$meta.Activate() | Out-Null
$meta.Range("A1").Select() | Out-Null
